# Auto-generated: apply scheduled-runner market-price refresh to Leve profit tables
$wb = $excel.ActiveWorkbook

# ----- Sheet ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 823230.5600000001
$ws.Range("J17").Value = 823230.5600000001
$ws.Range("L17").Value = 2469691.68
$ws.Range("N17").Value = -2470027.68
$ws.Range("H18").Value = 2924.75
$ws.Range("I18").Value = 2109.7
$ws.Range("K18").Value = 2109.7
$ws.Range("M18").Value = -1825.7
$ws.Range("H127").Value = 1536810.8
$ws.Range("I127").Value = 493.5
$ws.Range("J127").Value = 2151337.5
$ws.Range("K127").Value = 1480.5
$ws.Range("L127").Value = 6454012.5
$ws.Range("M127").Value = 3479.5
$ws.Range("N127").Value = -6463932.5
$ws.Range("H132").Value = 4635544
$ws.Range("I132").Value = 5449.294
$ws.Range("J132").Value = 12506705
$ws.Range("K132").Value = 16347.882
$ws.Range("L132").Value = 37520115
$ws.Range("M132").Value = -13817.882
$ws.Range("N132").Value = -37525175
$ws.Range("H137").Value = 5130936.5
$ws.Range("I137").Value = 1460.96
$ws.Range("J137").Value = 14290714
$ws.Range("K137").Value = 4382.88
$ws.Range("L137").Value = 42872142
$ws.Range("M137").Value = -1832.88
$ws.Range("N137").Value = -42877242
$ws.Range("H138").Value = 4718756.5
$ws.Range("I138").Value = 1580.3784
$ws.Range("J138").Value = 15627226
$ws.Range("K138").Value = 4741.135200000001
$ws.Range("L138").Value = 46881678
$ws.Range("M138").Value = 398.8647999999994
$ws.Range("N138").Value = -46891958

# ----- Sheet ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1213.2
$ws.Range("I2").Value = 852.1
$ws.Range("J2").Value = 1453.9333
$ws.Range("K2").Value = 852.1
$ws.Range("L2").Value = 1453.9333
$ws.Range("M2").Value = -739.1
$ws.Range("N2").Value = -1679.9333
$ws.Range("H32").Value = 9634.546
$ws.Range("I32").Value = 11999.081
$ws.Range("J32").Value = 4774.1113
$ws.Range("K32").Value = 11999.081
$ws.Range("L32").Value = 4774.1113
$ws.Range("M32").Value = -11712.081
$ws.Range("N32").Value = -5348.1113
$ws.Range("H45").Value = 3242.6667
$ws.Range("I45").Value = 3702
$ws.Range("J45").Value = 2783.3333
$ws.Range("K45").Value = 3702
$ws.Range("L45").Value = 2783.3333
$ws.Range("M45").Value = -3325
$ws.Range("N45").Value = -3537.3333
$ws.Range("H74").Value = 11630464
$ws.Range("I74").Value = 15626354
$ws.Range("J74").Value = 6058.364
$ws.Range("K74").Value = 15626354
$ws.Range("L74").Value = 6058.364
$ws.Range("M74").Value = -15625480
$ws.Range("N74").Value = -7806.364
$ws.Range("H77").Value = 11630464
$ws.Range("I77").Value = 15626354
$ws.Range("J77").Value = 6058.364
$ws.Range("K77").Value = 78131770
$ws.Range("L77").Value = 30291.82
$ws.Range("M77").Value = -78127402
$ws.Range("N77").Value = -39027.82
$ws.Range("H116").Value = 1213.2
$ws.Range("I116").Value = 852.1
$ws.Range("J116").Value = 1453.9333
$ws.Range("K116").Value = 852.1
$ws.Range("L116").Value = 1453.9333
$ws.Range("M116").Value = 1441.9
$ws.Range("N116").Value = -6041.9333
$ws.Range("H122").Value = 5283.483
$ws.Range("I122").Value = 6746.65
$ws.Range("J122").Value = 2032
$ws.Range("K122").Value = 20239.95
$ws.Range("L122").Value = 6096
$ws.Range("M122").Value = -17789.95
$ws.Range("N122").Value = -10996
$ws.Range("H132").Value = 4099868.5
$ws.Range("I132").Value = 5320173
$ws.Range("J132").Value = 3131.0715
$ws.Range("K132").Value = 15960519
$ws.Range("L132").Value = 9393.2145
$ws.Range("M132").Value = -15957989
$ws.Range("N132").Value = -14453.2145

# ----- Sheet BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1213.2
$ws.Range("I3").Value = 852.1
$ws.Range("J3").Value = 1453.9333
$ws.Range("K3").Value = 852.1
$ws.Range("L3").Value = 1453.9333
$ws.Range("M3").Value = -738.1
$ws.Range("N3").Value = -1681.9333

# ----- Sheet CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 976.6667
$ws.Range("I16").Value = 750.1111
$ws.Range("J16").Value = 1203.2222
$ws.Range("K16").Value = 750.1111
$ws.Range("L16").Value = 1203.2222
$ws.Range("M16").Value = -463.1111
$ws.Range("N16").Value = -1777.2222
$ws.Range("H31").Value = 5053735.5
$ws.Range("I31").Value = 4990
$ws.Range("J31").Value = 10753932
$ws.Range("K31").Value = 4990
$ws.Range("L31").Value = 10753932
$ws.Range("M31").Value = -4695
$ws.Range("N31").Value = -10754522
$ws.Range("H34").Value = 5053735.5
$ws.Range("I34").Value = 4990
$ws.Range("J34").Value = 10753932
$ws.Range("K34").Value = 4990
$ws.Range("L34").Value = 10753932
$ws.Range("M34").Value = -4788
$ws.Range("N34").Value = -10754336
$ws.Range("H113").Value = 976.6667
$ws.Range("I113").Value = 750.1111
$ws.Range("J113").Value = 1203.2222
$ws.Range("K113").Value = 750.1111
$ws.Range("L113").Value = 1203.2222
$ws.Range("M113").Value = 1419.8889
$ws.Range("N113").Value = -5543.2222

# ----- Sheet CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 721582.6
$ws.Range("J33").Value = 88.8
$ws.Range("L33").Value = 532.8
$ws.Range("N33").Value = -1098.8
$ws.Range("H99").Value = 2008.3334
$ws.Range("I99").Value = 1525
$ws.Range("J99").Value = 2250
$ws.Range("K99").Value = 4575
$ws.Range("L99").Value = 6750
$ws.Range("M99").Value = -2329
$ws.Range("N99").Value = -11242
$ws.Range("H122").Value = 966.05
$ws.Range("J122").Value = 604.7778
$ws.Range("L122").Value = 5443.000199999999
$ws.Range("N122").Value = -10343.0002
$ws.Range("H134").Value = 3861.5386
$ws.Range("I134").Value = 2534.375
$ws.Range("J134").Value = 5985
$ws.Range("K134").Value = 7603.125
$ws.Range("L134").Value = 17955
$ws.Range("M134").Value = -2533.125
$ws.Range("N134").Value = -28095

# ----- Sheet GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H88").Value = 36500
$ws.Range("J88").Value = 36500
$ws.Range("L88").Value = 36500
$ws.Range("N88").Value = -37402
$ws.Range("H91").Value = 36500
$ws.Range("J91").Value = 36500
$ws.Range("L91").Value = 36500
$ws.Range("N91").Value = -39620
$ws.Range("H102").Value = 2950.92
$ws.Range("I102").Value = 3526.7778
$ws.Range("J102").Value = 1470.1428
$ws.Range("K102").Value = 3526.7778
$ws.Range("L102").Value = 1470.1428
$ws.Range("M102").Value = -1904.7778
$ws.Range("N102").Value = -4714.1428
$ws.Range("H126").Value = 4651.6
$ws.Range("I126").Value = 2862
$ws.Range("J126").Value = 5615.231
$ws.Range("K126").Value = 8586
$ws.Range("L126").Value = 16845.693
$ws.Range("M126").Value = -6116
$ws.Range("N126").Value = -21785.693
$ws.Range("H132").Value = 2778.0188
$ws.Range("I132").Value = 1759.7297
$ws.Range("J132").Value = 5132.8125
$ws.Range("K132").Value = 5279.189100000001
$ws.Range("L132").Value = 15398.4375
$ws.Range("M132").Value = -2749.189100000001
$ws.Range("N132").Value = -20458.4375

# ----- Sheet LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4472.3955
$ws.Range("I7").Value = 3889.5715
$ws.Range("J7").Value = 5560.3335
$ws.Range("K7").Value = 3889.5715
$ws.Range("L7").Value = 5560.3335
$ws.Range("M7").Value = -3777.5715
$ws.Range("N7").Value = -5784.3335
$ws.Range("H40").Value = 4589.2583
$ws.Range("I40").Value = 4574.35
$ws.Range("J40").Value = 4616.364
$ws.Range("K40").Value = 4574.35
$ws.Range("L40").Value = 4616.364
$ws.Range("M40").Value = -4438.35
$ws.Range("N40").Value = -4888.364
$ws.Range("H55").Value = 290.875
$ws.Range("I55").Value = 143.54546
$ws.Range("J55").Value = 415.53845
$ws.Range("K55").Value = 143.54546
$ws.Range("L55").Value = 415.53845
$ws.Range("M55").Value = 29.45454000000001
$ws.Range("N55").Value = -761.53845
$ws.Range("H126").Value = 4472.3955
$ws.Range("I126").Value = 3889.5715
$ws.Range("J126").Value = 5560.3335
$ws.Range("K126").Value = 11668.7145
$ws.Range("L126").Value = 16681.0005
$ws.Range("M126").Value = -9198.7145
$ws.Range("N126").Value = -21621.0005

# ----- Sheet WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H132").Value = 4285.892
$ws.Range("I132").Value = 5730.3335
$ws.Range("J132").Value = 1619.2307
$ws.Range("K132").Value = 17191.0005
$ws.Range("L132").Value = 4857.6921
$ws.Range("M132").Value = -14661.0005
$ws.Range("N132").Value = -9917.6921
$ws.Range("H136").Value = 915.1515000000001
$ws.Range("J136").Value = 1080
$ws.Range("L136").Value = 3240
$ws.Range("N136").Value = -8340
$ws.Range("H138").Value = 60407.332
$ws.Range("J138").Value = 60407.332
$ws.Range("L138").Value = 60407.332
$ws.Range("N138").Value = -70687.33199999999
